$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 54.857143
$ws.Range("I6").Value = 61.17647
$ws.Range("K6").Value = 183.52941
$ws.Range("M6").Value = -71.52941000000001
$ws.Range("H9").Value = 3675.6924
$ws.Range("I9").Value = 3642.7778
$ws.Range("J9").Value = 3749.75
$ws.Range("K9").Value = 3642.7778
$ws.Range("L9").Value = 3749.75
$ws.Range("M9").Value = -3473.7778
$ws.Range("N9").Value = -4087.75
$ws.Range("H38").Value = 2696.0908
$ws.Range("J38").Value = 3631.125
$ws.Range("L38").Value = 10893.375
$ws.Range("N38").Value = -11637.375
$ws.Range("H54").Value = 265000
$ws.Range("I54").Value = 265000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 265000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -264514
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 1427.6428
$ws.Range("J58").Value = 2626.8572
$ws.Range("L58").Value = 7880.571599999999
$ws.Range("N58").Value = -8180.571599999999
$ws.Range("H99").Value = 1585.0834
$ws.Range("I99").Value = 1593.2727
$ws.Range("K99").Value = 4779.8181
$ws.Range("M99").Value = -3281.8181
$ws.Range("H101").Value = 1435.7
$ws.Range("I101").Value = 1665.2858
$ws.Range("J101").Value = 900
$ws.Range("K101").Value = 4995.857400000001
$ws.Range("L101").Value = 2700
$ws.Range("M101").Value = -3373.857400000001
$ws.Range("N101").Value = -5944
$ws.Range("H103").Value = 366.33334
$ws.Range("I103").Value = 299
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 897
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = -311
$ws.Range("N103").Value = -2372
$ws.Range("H106").Value = 13782.889
$ws.Range("I106").Value = 16649.715
$ws.Range("K106").Value = 16649.715
$ws.Range("M106").Value = -16018.715
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H112").Value = 3818.7778
$ws.Range("J112").Value = 4484.143
$ws.Range("L112").Value = 13452.429
$ws.Range("N112").Value = -15668.429
$ws.Range("H116").Value = 2180.95
$ws.Range("I116").Value = 2139.1753
$ws.Range("J116").Value = 3531.6667
$ws.Range("K116").Value = 2139.1753
$ws.Range("L116").Value = 3531.6667
$ws.Range("M116").Value = 1302.8247
$ws.Range("N116").Value = -10415.6667
$ws.Range("H132").Value = 7939.2085
$ws.Range("J132").Value = 16355
$ws.Range("L132").Value = 49065
$ws.Range("N132").Value = -54125
$ws.Range("H137").Value = 7192.4507
$ws.Range("I137").Value = 2337.4546
$ws.Range("K137").Value = 7012.3638
$ws.Range("M137").Value = -4462.3638
$ws.Range("H138").Value = 3555.4084
$ws.Range("I138").Value = 3727.4285
$ws.Range("J138").Value = 3443.3953
$ws.Range("K138").Value = 11182.2855
$ws.Range("L138").Value = 10330.1859
$ws.Range("M138").Value = -6042.2855
$ws.Range("N138").Value = -20610.1859
$ws.Range("H141").Value = 5664.8696
$ws.Range("I141").Value = 5238.5557
$ws.Range("K141").Value = 15715.6671
$ws.Range("M141").Value = -10535.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22296.537
$ws.Range("I32").Value = 4759.508
$ws.Range("J32").Value = 56822.562
$ws.Range("K32").Value = 4759.508
$ws.Range("L32").Value = 56822.562
$ws.Range("M32").Value = -4472.508
$ws.Range("N32").Value = -57396.562
$ws.Range("H61").Value = 17509.92
$ws.Range("I61").Value = 10162.154
$ws.Range("K61").Value = 10162.154
$ws.Range("M61").Value = -9950.154
$ws.Range("H74").Value = 16204.467
$ws.Range("I74").Value = 2804.4119
$ws.Range("K74").Value = 2804.4119
$ws.Range("M74").Value = -1930.4119
$ws.Range("H77").Value = 16204.467
$ws.Range("I77").Value = 2804.4119
$ws.Range("K77").Value = 14022.0595
$ws.Range("M77").Value = -9654.059499999999
$ws.Range("H97").Value = 1568.6
$ws.Range("I97").Value = 1795.2
$ws.Range("J97").Value = 1342
$ws.Range("K97").Value = 1795.2
$ws.Range("L97").Value = 1342
$ws.Range("M97").Value = -1299.2
$ws.Range("N97").Value = -2334
$ws.Range("H136").Value = 17509.92
$ws.Range("I136").Value = 10162.154
$ws.Range("K136").Value = 30486.462
$ws.Range("M136").Value = -27936.462
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6382.077
$ws.Range("I86").Value = 5797.5
$ws.Range("K86").Value = 5797.5
$ws.Range("M86").Value = -4674.5
$ws.Range("H89").Value = 6382.077
$ws.Range("I89").Value = 5797.5
$ws.Range("K89").Value = 28987.5
$ws.Range("M89").Value = -23371.5
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H94").Value = 3541.182
$ws.Range("I94").Value = 1994.8889
$ws.Range("J94").Value = 10499.5
$ws.Range("K94").Value = 1994.8889
$ws.Range("L94").Value = 10499.5
$ws.Range("M94").Value = -1543.8889
$ws.Range("N94").Value = -11401.5
$ws.Range("H96").Value = 29865.584
$ws.Range("I96").Value = 9732
$ws.Range("K96").Value = 9732
$ws.Range("M96").Value = -6986
$ws.Range("H134").Value = 8436.366
$ws.Range("I134").Value = 2352.6086
$ws.Range("J134").Value = 16210.056
$ws.Range("K134").Value = 7057.825800000001
$ws.Range("L134").Value = 48630.16800000001
$ws.Range("M134").Value = -4522.825800000001
$ws.Range("N134").Value = -53700.16800000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4149.8
$ws.Range("I22").Value = 6000.3335
$ws.Range("K22").Value = 6000.3335
$ws.Range("M22").Value = -5650.3335
$ws.Range("H58").Value = 17876.842
$ws.Range("I58").Value = 11603.667
$ws.Range("J58").Value = 20772.154
$ws.Range("K58").Value = 11603.667
$ws.Range("L58").Value = 20772.154
$ws.Range("M58").Value = -11400.667
$ws.Range("N58").Value = -21178.154
$ws.Range("H59").Value = 35529.7
$ws.Range("J59").Value = 39649.125
$ws.Range("L59").Value = 39649.125
$ws.Range("N59").Value = -41939.125
$ws.Range("H107").Value = 1329.3334
$ws.Range("I107").Value = 1008.25
$ws.Range("K107").Value = 1008.25
$ws.Range("M107").Value = 911.75
$ws.Range("H132").Value = 5433.1284
$ws.Range("I132").Value = 1609.8518
$ws.Range("K132").Value = 4829.555399999999
$ws.Range("M132").Value = -2299.555399999999
$ws.Range("H136").Value = 17876.842
$ws.Range("I136").Value = 11603.667
$ws.Range("J136").Value = 20772.154
$ws.Range("K136").Value = 34811.001
$ws.Range("L136").Value = 62316.462
$ws.Range("M136").Value = -32261.001
$ws.Range("N136").Value = -67416.462
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1591740.4
$ws.Range("J41").Value = 1819103.2
$ws.Range("L41").Value = 5457309.6
$ws.Range("N41").Value = -5457985.6
$ws.Range("H104").Value = 4117559.8
$ws.Range("I104").Value = 5000
$ws.Range("K104").Value = 15000
$ws.Range("M104").Value = -12379
$ws.Range("H131").Value = 1471.91
$ws.Range("J131").Value = 1476.7755
$ws.Range("L131").Value = 4430.3265
$ws.Range("N131").Value = -14510.3265
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19356.715
$ws.Range("J80").Value = 30199
$ws.Range("L80").Value = 30199
$ws.Range("N80").Value = -32195
$ws.Range("H83").Value = 19356.715
$ws.Range("J83").Value = 30199
$ws.Range("L83").Value = 150995
$ws.Range("N83").Value = -160979
$ws.Range("H103").Value = 140000
$ws.Range("J103").Value = 140000
$ws.Range("L103").Value = 140000
$ws.Range("N103").Value = -142344
$ws.Range("H107").Value = 1480.0588
$ws.Range("I107").Value = 1347.1
$ws.Range("K107").Value = 1347.1
$ws.Range("M107").Value = 572.9000000000001
$ws.Range("H132").Value = 7606.05
$ws.Range("I132").Value = 2881.75
$ws.Range("K132").Value = 8645.25
$ws.Range("M132").Value = -6115.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1606.2858
$ws.Range("I16").Value = 1774.5217
$ws.Range("K16").Value = 1774.5217
$ws.Range("M16").Value = -1604.5217
$ws.Range("H40").Value = 13000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 2910.2964
$ws.Range("I46").Value = 1353.8462
$ws.Range("J46").Value = 4355.5713
$ws.Range("K46").Value = 1353.8462
$ws.Range("L46").Value = 4355.5713
$ws.Range("M46").Value = -1165.8462
$ws.Range("N46").Value = -4731.5713
$ws.Range("H122").Value = 6722.222
$ws.Range("I122").Value = 6357.143
$ws.Range("K122").Value = 19071.429
$ws.Range("M122").Value = -16621.429
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 75000
$ws.Range("J56").Value = 75000
$ws.Range("L56").Value = 75000
$ws.Range("N56").Value = -76428
